# Refresh NATMI ligand-receptor (Fn1-Itga2) edge-weight statistics with updated TPM values.
# Columns G-J: ligand stats (per "Sending cluster"); K-P: receptor stats (per "Target cluster");
# Q-T: derived edge weights/specificities recomputed from the refreshed ligand & receptor stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 32.93949833333333
$ws.Range("H2").Value = 98.818495
$ws.Range("I2").Value = 0.02571831923682078
$ws.Range("J2").Value = 0.02571831923682077
$ws.Range("M2").Value = 3.339352
$ws.Range("N2").Value = 10.018056
$ws.Range("O2").Value = 0.6054960700393903
$ws.Range("P2").Value = 0.6054960700393903
$ws.Range("Q2").Value = 109.9965796384133
$ws.Range("R2").Value = 989.9692167457201
$ws.Range("S2").Value = 0.01557234122591343
$ws.Range("T2").Value = 0.01557234122591343

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 32.93949833333333
$ws.Range("H3").Value = 98.818495
$ws.Range("I3").Value = 0.02571831923682078
$ws.Range("J3").Value = 0.02571831923682077
$ws.Range("O3").Value = 0.2540955070726236
$ws.Range("P3").Value = 0.2540955070726236
$ws.Range("Q3").Value = 46.15989774741001
$ws.Range("R3").Value = 415.43907972669
$ws.Range("S3").Value = 0.006534909367535586
$ws.Range("T3").Value = 0.006534909367535584

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("G4").Value = 32.93949833333333
$ws.Range("H4").Value = 98.818495
$ws.Range("I4").Value = 0.02571831923682078
$ws.Range("J4").Value = 0.02571831923682077
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1338136666666667
$ws.Range("N4").Value = 0.401441
$ws.Range("O4").Value = 0.02426328499787613
$ws.Range("P4").Value = 0.02426328499787612
$ws.Range("Q4").Value = 4.40775505014389
$ws.Range("R4").Value = 39.669795451295
$ws.Range("S4").Value = 0.0006240109093093425
$ws.Range("T4").Value = 0.0006240109093093424

# Row 5: ECs -> MuSCs
$ws.Range("G5").Value = 32.93949833333333
$ws.Range("H5").Value = 98.818495
$ws.Range("I5").Value = 0.02571831923682078
$ws.Range("J5").Value = 0.02571831923682077
$ws.Range("M5").Value = 0.6405483333333334
$ws.Range("N5").Value = 1.921645
$ws.Range("O5").Value = 0.11614513789011
$ws.Range("P5").Value = 0.11614513789011
$ws.Range("Q5").Value = 21.09934075825278
$ws.Range("R5").Value = 189.894066824275
$ws.Range("S5").Value = 0.002987057734062419
$ws.Range("T5").Value = 0.002987057734062418

# Row 6: FAPs -> ECs
$ws.Range("I6").Value = 0.4140443484779395
$ws.Range("J6").Value = 0.4140443484779395
$ws.Range("M6").Value = 3.339352
$ws.Range("N6").Value = 10.018056
$ws.Range("O6").Value = 0.6054960700393903
$ws.Range("P6").Value = 0.6054960700393903
$ws.Range("Q6").Value = 1770.856864004717
$ws.Range("R6").Value = 15937.71177604246
$ws.Range("S6").Value = 0.2507022258254122
$ws.Range("T6").Value = 0.2507022258254122

# Row 7: FAPs -> FAPs
$ws.Range("I7").Value = 0.4140443484779395
$ws.Range("J7").Value = 0.4140443484779395
$ws.Range("O7").Value = 0.2540955070726236
$ws.Range("P7").Value = 0.2540955070726236
$ws.Range("S7").Value = 0.1052068086770561
$ws.Range("T7").Value = 0.1052068086770561

# Row 8: FAPs -> Inflammatory-Mac
$ws.Range("I8").Value = 0.4140443484779395
$ws.Range("J8").Value = 0.4140443484779395
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1338136666666667
$ws.Range("N8").Value = 0.401441
$ws.Range("O8").Value = 0.02426328499787613
$ws.Range("P8").Value = 0.02426328499787612
$ws.Range("Q8").Value = 70.9613272617879
$ws.Range("R8").Value = 638.6519453560911
$ws.Range("S8").Value = 0.01004607602888018
$ws.Range("T8").Value = 0.01004607602888018

# Row 9: FAPs -> MuSCs
$ws.Range("I9").Value = 0.4140443484779395
$ws.Range("J9").Value = 0.4140443484779395
$ws.Range("M9").Value = 0.6405483333333334
$ws.Range("N9").Value = 1.921645
$ws.Range("O9").Value = 0.11614513789011
$ws.Range("P9").Value = 0.11614513789011
$ws.Range("Q9").Value = 339.6824931334328
$ws.Range("R9").Value = 3057.142438200895
$ws.Range("S9").Value = 0.04808923794659105
$ws.Range("T9").Value = 0.04808923794659105

# Row 10: Inflammatory-Mac -> ECs
$ws.Range("G10").Value = 422.1807963333333
$ws.Range("H10").Value = 1266.542389
$ws.Range("I10").Value = 0.3296279860087694
$ws.Range("J10").Value = 0.3296279860087693
$ws.Range("M10").Value = 3.339352
$ws.Range("N10").Value = 10.018056
$ws.Range("O10").Value = 0.6054960700393903
$ws.Range("P10").Value = 0.6054960700393903
$ws.Range("Q10").Value = 1409.810286597309
$ws.Range("R10").Value = 12688.29257937579
$ws.Range("S10").Value = 0.199588450103309
$ws.Range("T10").Value = 0.1995884501033089

# Row 11: Inflammatory-Mac -> FAPs
$ws.Range("G11").Value = 422.1807963333333
$ws.Range("H11").Value = 1266.542389
$ws.Range("I11").Value = 0.3296279860087694
$ws.Range("J11").Value = 0.3296279860087693
$ws.Range("O11").Value = 0.2540955070726236
$ws.Range("P11").Value = 0.2540955070726236
$ws.Range("Q11").Value = 591.6247476649021
$ws.Range("R11").Value = 5324.622728984118
$ws.Range("S11").Value = 0.08375699025022594
$ws.Range("T11").Value = 0.0837569902502259

# Row 12: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("G12").Value = 422.1807963333333
$ws.Range("H12").Value = 1266.542389
$ws.Range("I12").Value = 0.3296279860087694
$ws.Range("J12").Value = 0.3296279860087693
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1338136666666667
$ws.Range("N12").Value = 0.401441
$ws.Range("O12").Value = 0.02426328499787613
$ws.Range("P12").Value = 0.02426328499787612
$ws.Range("Q12").Value = 56.49356035361657
$ws.Range("R12").Value = 508.4420431825491
$ws.Range("S12").Value = 0.007997857767806696
$ws.Range("T12").Value = 0.007997857767806692

# Row 13: Inflammatory-Mac -> MuSCs
$ws.Range("G13").Value = 422.1807963333333
$ws.Range("H13").Value = 1266.542389
$ws.Range("I13").Value = 0.3296279860087694
$ws.Range("J13").Value = 0.3296279860087693
$ws.Range("M13").Value = 0.6405483333333334
$ws.Range("N13").Value = 1.921645
$ws.Range("O13").Value = 0.11614513789011
$ws.Range("P13").Value = 0.11614513789011
$ws.Range("Q13").Value = 270.4272054566561
$ws.Range("R13").Value = 2433.844849109905
$ws.Range("S13").Value = 0.03828468788742778
$ws.Range("T13").Value = 0.03828468788742777

# Row 14: MuSCs -> ECs
$ws.Range("G14").Value = 16.509264
$ws.Range("H14").Value = 49.527792
$ws.Range("I14").Value = 0.01289001179132366
$ws.Range("J14").Value = 0.01289001179132366
$ws.Range("M14").Value = 3.339352
$ws.Range("N14").Value = 10.018056
$ws.Range("O14").Value = 0.6054960700393903
$ws.Range("P14").Value = 0.6054960700393903
$ws.Range("Q14").Value = 55.130243756928
$ws.Range("R14").Value = 496.1721938123521
$ws.Range("S14").Value = 0.007804851482407877
$ws.Range("T14").Value = 0.007804851482407877

# Row 15: MuSCs -> FAPs
$ws.Range("G15").Value = 16.509264
$ws.Range("H15").Value = 49.527792
$ws.Range("I15").Value = 0.01289001179132366
$ws.Range("J15").Value = 0.01289001179132366
$ws.Range("O15").Value = 0.2540955070726236
$ws.Range("P15").Value = 0.2540955070726236
$ws.Range("Q15").Value = 23.135323143456
$ws.Range("R15").Value = 208.217908291104
$ws.Range("S15").Value = 0.003275294082288483
$ws.Range("T15").Value = 0.003275294082288482

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Range("G16").Value = 16.509264
$ws.Range("H16").Value = 49.527792
$ws.Range("I16").Value = 0.01289001179132366
$ws.Range("J16").Value = 0.01289001179132366
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1338136666666667
$ws.Range("N16").Value = 0.401441
$ws.Range("O16").Value = 0.02426328499787613
$ws.Range("P16").Value = 0.02426328499787612
$ws.Range("Q16").Value = 2.209165149808
$ws.Range("R16").Value = 19.882486348272
$ws.Range("S16").Value = 0.0003127540297188697
$ws.Range("T16").Value = 0.0003127540297188697

# Row 17: MuSCs -> MuSCs
$ws.Range("G17").Value = 16.509264
$ws.Range("H17").Value = 49.527792
$ws.Range("I17").Value = 0.01289001179132366
$ws.Range("J17").Value = 0.01289001179132366
$ws.Range("M17").Value = 0.6405483333333334
$ws.Range("N17").Value = 1.921645
$ws.Range("O17").Value = 0.11614513789011
$ws.Range("P17").Value = 0.11614513789011
$ws.Range("Q17").Value = 10.57498153976
$ws.Range("R17").Value = 95.17483385784
$ws.Range("S17").Value = 0.001497112196908431
$ws.Range("T17").Value = 0.00149711219690843

# Row 18: Neutrophils -> ECs
$ws.Range("G18").Value = 44.62094166666666
$ws.Range("H18").Value = 133.862825
$ws.Range("I18").Value = 0.03483889192294087
$ws.Range("J18").Value = 0.03483889192294087
$ws.Range("M18").Value = 3.339352
$ws.Range("N18").Value = 10.018056
$ws.Range("O18").Value = 0.6054960700393903
$ws.Range("P18").Value = 0.6054960700393903
$ws.Range("Q18").Value = 149.0050307964667
$ws.Range("R18").Value = 1341.0452771682
$ws.Range("S18").Value = 0.02109481214386776
$ws.Range("T18").Value = 0.02109481214386776

# Row 19: Neutrophils -> FAPs
$ws.Range("G19").Value = 44.62094166666666
$ws.Range("H19").Value = 133.862825
$ws.Range("I19").Value = 0.03483889192294087
$ws.Range("J19").Value = 0.03483889192294087
$ws.Range("O19").Value = 0.2540955070726236
$ws.Range("P19").Value = 0.2540955070726236
$ws.Range("Q19").Value = 62.52973508835
$ws.Range("R19").Value = 562.76761579515
$ws.Range("S19").Value = 0.008852405909007993
$ws.Range("T19").Value = 0.008852405909007991

# Row 20: Neutrophils -> Inflammatory-Mac
$ws.Range("G20").Value = 44.62094166666666
$ws.Range("H20").Value = 133.862825
$ws.Range("I20").Value = 0.03483889192294087
$ws.Range("J20").Value = 0.03483889192294087
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.1338136666666667
$ws.Range("N20").Value = 0.401441
$ws.Range("O20").Value = 0.02426328499787613
$ws.Range("P20").Value = 0.02426328499787612
$ws.Range("Q20").Value = 5.970891814536111
$ws.Range("R20").Value = 53.738026330825
$ws.Range("S20").Value = 0.0008453059637365189
$ws.Range("T20").Value = 0.0008453059637365188

# Row 21: Neutrophils -> MuSCs
$ws.Range("G21").Value = 44.62094166666666
$ws.Range("H21").Value = 133.862825
$ws.Range("I21").Value = 0.03483889192294087
$ws.Range("J21").Value = 0.03483889192294087
$ws.Range("M21").Value = 0.6405483333333334
$ws.Range("N21").Value = 1.921645
$ws.Range("O21").Value = 0.11614513789011
$ws.Range("P21").Value = 0.11614513789011
$ws.Range("Q21").Value = 28.58186981634722
$ws.Range("R21").Value = 257.236828347125
$ws.Range("S21").Value = 0.004046367906328609
$ws.Range("T21").Value = 0.004046367906328608

# Row 22: Resolving-Mac -> ECs
$ws.Range("G22").Value = 234.229538
$ws.Range("H22").Value = 702.6886139999999
$ws.Range("I22").Value = 0.1828804425622059
$ws.Range("J22").Value = 0.1828804425622059
$ws.Range("M22").Value = 3.339352
$ws.Range("N22").Value = 10.018056
$ws.Range("O22").Value = 0.6054960700393903
$ws.Range("P22").Value = 0.6054960700393903
$ws.Range("Q22").Value = 782.1748761793759
$ws.Range("R22").Value = 7039.573885614384
$ws.Range("S22").Value = 0.1107333892584801
$ws.Range("T22").Value = 0.1107333892584801

# Row 23: Resolving-Mac -> FAPs
$ws.Range("G23").Value = 234.229538
$ws.Range("H23").Value = 702.6886139999999
$ws.Range("I23").Value = 0.1828804425622059
$ws.Range("J23").Value = 0.1828804425622059
$ws.Range("O23").Value = 0.2540955070726236
$ws.Range("P23").Value = 0.2540955070726236
$ws.Range("Q23").Value = 328.238499994452
$ws.Range("R23").Value = 2954.146499950068
$ws.Range("S23").Value = 0.04646909878650952
$ws.Range("T23").Value = 0.04646909878650951

# Row 24: Resolving-Mac -> Inflammatory-Mac
$ws.Range("G24").Value = 234.229538
$ws.Range("H24").Value = 702.6886139999999
$ws.Range("I24").Value = 0.1828804425622059
$ws.Range("J24").Value = 0.1828804425622059
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.1338136666666667
$ws.Range("N24").Value = 0.401441
$ws.Range("O24").Value = 0.02426328499787613
$ws.Range("P24").Value = 0.02426328499787612
$ws.Range("Q24").Value = 31.34311332141933
$ws.Range("R24").Value = 282.088019892774
$ws.Range("S24").Value = 0.004437280298424516
$ws.Range("T24").Value = 0.004437280298424515

# Row 25: Resolving-Mac -> MuSCs
$ws.Range("G25").Value = 234.229538
$ws.Range("H25").Value = 702.6886139999999
$ws.Range("I25").Value = 0.1828804425622059
$ws.Range("J25").Value = 0.1828804425622059
$ws.Range("M25").Value = 0.6405483333333334
$ws.Range("N25").Value = 1.921645
$ws.Range("O25").Value = 0.11614513789011
$ws.Range("P25").Value = 0.11614513789011
$ws.Range("Q25").Value = 150.0353401833366
$ws.Range("R25").Value = 1350.31806165003
$ws.Range("S25").Value = 0.02124067421879175
$ws.Range("T25").Value = 0.02124067421879175
